$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[ Thanos%Karatzias%null%1,                            Elsie%Yan%null%1,                            Sally%Jowett%null%1,                         Thanos%Karatzias%null%1,                         Elsie%Yan%null%1,                         Sally%Jowett%null%1]"
$ws.Range("E3").Value = "[Wayne R.%Lawrence%xref no email%1,       Ziqiang%Lin%xref no email%1,       Emily A.%Lipton%xref no email%1,       Guthrie%Birkhead%xref no email%1,       Michael%Primeau%xref no email%1,       Guang-Hui%Dong%xref no email%1,       Shao%Lin%xref no email%1]"
$ws.Range("E4").Value = "[Alain%Vanasse%NULL%1,                           Alan%Cohen%NULL%1,                           Josiane%Courteau%NULL%1,                           Patrick%Bergeron%NULL%1,                           Roxanne%Dault%NULL%1,                           Pierre%Gosselin%NULL%1,                           Claudia%Blais%NULL%1,                           Diane%Bélanger%NULL%1,                           Louis%Rochette%NULL%1,                           Fateh%Chebana%NULL%1,                           Jan%Semenza%NULL%2,                           Jan%Semenza%NULL%0]"
$ws.Range("E5").Value = "[ Leonardo%Trasande%null%1,                            Tony T.%Koshy%null%1,                            Joseph%Gilbert%null%1,                            Lauren K.%Burdine%null%1,                            Michael%Marmor%null%1,                            Xiaoxia%Han%null%1,                            Yongzhao%Shao%null%1,                            Claude%Chemtob%null%1,                            Teresa M.%Attina%null%1,                            Elaine M.%Urbina%null%1,                         Leonardo%Trasande%null%1,                         Tony T.%Koshy%null%1,                         Joseph%Gilbert%null%1,                         Lauren K.%Burdine%null%1,                         Michael%Marmor%null%1,                         Xiaoxia%Han%null%1,                         Yongzhao%Shao%null%1,                         Claude%Chemtob%null%1,                         Teresa M.%Attina%null%1,                         Elaine M.%Urbina%null%1]"
$ws.Range("E6").Value = "[Ninon A.%Becquart%NULL%1,                           Elena N.%Naumova%NULL%1,                           Gitanjali%Singh%NULL%1,                           Kenneth K. H.%Chui%NULL%1]"
$ws.Range("E7").Value = "[Hannah T.%Jordan%NULL%1,                           Steven D.%Stellman%NULL%1,                           Alfredo%Morabia%NULL%1,                           Sara A.%Miller‐Archie%NULL%1,                           Howard%Alper%NULL%1,                           Zoey%Laskaris%NULL%1,                           Robert M.%Brackbill%NULL%1,                           James E.%Cone%NULL%1]"
$ws.Range("E8").Value = "[So Yeon%Kong%xref no email%1,       Kyoung Jun%Song%xref no email%1,       Sang Do%Shin%xref no email%1,       Young Sun%Ro%xref no email%1]"
$ws.Range("E9").Value = "[Hector%F. Myers%coreGivesNoEmail%1,                          Karen%A. Matthews%coreGivesNoEmail%1,                          Nataria%T. Joseph%coreGivesNoEmail%1]"
$ws.Range("E10").Value = "[ Nathan%McKinney%null%1,                  Chris%Houser%null%1,                  Klaus%Meyer-Arendt%null%1]"
$ws.Range("E11").Value = "[Anja J. E.%Dirkzwager%xref no email%1,       Peter G.%van der Velden%xref no email%1,       Linda%Grievink%xref no email%1,       C Joris%Yzermans%xref no email%1]"
$ws.Range("E12").Value = "[ Zhen%Jiao%null%1,                            Socrates V.%Kakoulides%null%1,                            John%Moscona%null%1,                            Jabar%Whittier%null%1,                            Sudesh%Srivastav%null%1,                            Patrice%Delafontaine%null%2,                            Anand%Irimpen%null%1,                         Zhen%Jiao%null%1,                         Socrates V.%Kakoulides%null%1,                         John%Moscona%null%1,                         Jabar%Whittier%null%1,                         Sudesh%Srivastav%null%1,                         Patrice%Delafontaine%null%2,                         Anand%Irimpen%null%1]"
$ws.Range("E13").Value = "[Sandeep%Gautam%xref no email%1,       Jonathan%Menachem%xref no email%1,       Sudesh K.%Srivastav%xref no email%1,       Patrice%Delafontaine%xref no email%1,       Anand%Irimpen%xref no email%1]"
$ws.Range("E14").Value = "[ Yasuhiro%Nagayoshi%null%1,                            Shinya%Yumoto%null%1,                            Kazuhisa%Sakaguchi%null%1,                            Chiharu%Shudo%null%1,                            Shiro%Takino%null%1,                            Motohiro%Hashiyama%null%1,                            Yutaka%Kai%null%1,                            Yutaka%Kuroda%null%1,                            Hiroaki%Kawano%null%1,                            Hisao%Ogawa%null%1,                         Yasuhiro%Nagayoshi%null%1,                         Shinya%Yumoto%null%1,                         Kazuhisa%Sakaguchi%null%1,                         Chiharu%Shudo%null%1,                         Shiro%Takino%null%1,                         Motohiro%Hashiyama%null%1,                         Yutaka%Kai%null%1,                         Yutaka%Kuroda%null%1,                         Hiroaki%Kawano%null%1,                         Hisao%Ogawa%null%1]"
$ws.Range("E15").Value = "[ Hannah T.%Jordan%null%1,                            Sara A.%Miller-Archie%null%1,                            James E.%Cone%null%1,                            Alfredo%Morabia%null%1,                            Steven D.%Stellman%null%1,                         Hannah T.%Jordan%null%1,                         Sara A.%Miller-Archie%null%1,                         James E.%Cone%null%1,                         Alfredo%Morabia%null%1,                         Steven D.%Stellman%null%1]"
$ws.Range("E16").Value = "[Soyeon%Kim%xref no email%1,       Prathit A.%Kulkarni%xref no email%1,       Mangala%Rajan%xref no email%1,       Pauline%Thomas%xref no email%1,       Stella%Tsai%xref no email%1,       Christina%Tan%xref no email%1,       Amy%Davidow%xref no email%1]"
$ws.Range("E17").Value = "[Vivian A.%Fonseca%NULL%1,                           Hayden%Smith%NULL%1,                           Nitesh%Kuhadiya%NULL%1,                           Sharice M.%Leger%NULL%1,                           C. Lillian%Yau%NULL%1,                           Kristi%Reynolds%NULL%1,                           Lizheng%Shi%NULL%1,                           Roberta H.%McDuffie%NULL%1,                           Tina%Thethi%NULL%1,                           Jennifer%John-Kalarickal%NULL%1]"
$ws.Range("E18").Value = "[ Federico%Silva-Palacios%null%1,                            Ana Isabel%Casanegra%null%1,                            Alan%Shapiro%null%1,                            Minh%Phan%null%1,                            Beau%Hawkins%null%1,                            Ji%Li%null%1,                            Julie%Stoner%null%1,                            Alfonso%Tafur%null%1,                         Federico%Silva-Palacios%null%1,                         Ana Isabel%Casanegra%null%1,                         Alan%Shapiro%null%1,                         Minh%Phan%null%1,                         Beau%Hawkins%null%1,                         Ji%Li%null%1,                         Julie%Stoner%null%1,                         Alfonso%Tafur%null%1]"
$ws.Range("E19").Value = "[Tran Huu%Bich%NULL%1,                           La Ngoc%Quang%NULL%1,                           Le Thi%Thanh Ha%NULL%1,                           Tran Thi%Duc Hanh%NULL%1,                           Debarati%Guha-Sapir%NULL%1]"
$ws.Range("E20").Value = "[Bijlsma%J A%coreGivesNoEmail%1,                          Bouter%L M%coreGivesNoEmail%1,                          Bramsen%I%coreGivesNoEmail%1,                          Huizink%A C%coreGivesNoEmail%1,                          Slottje%P%coreGivesNoEmail%1,                          Smid%T%coreGivesNoEmail%1,                          Smidt%N%coreGivesNoEmail%1,                          Twisk%J W R%coreGivesNoEmail%1,                          van%der Ploeg H M%coreGivesNoEmail%1,                          van%Mechelen W%coreGivesNoEmail%1,                          Witteveen%A B%coreGivesNoEmail%1]"
$ws.Range("E21").Value = "[ Hannah T%Jordan%null%1,                            Robert M%Brackbill%null%1,                            James E%Cone%null%1,                            Indira%Debchoudhury%null%1,                            Mark R%Farfel%null%1,                            Carolyn M%Greene%null%1,                            James L%Hadler%null%1,                            Joseph%Kennedy%null%1,                            Jiehui%Li%null%1,                            Jonathan%Liff%null%1,                            Leslie%Stayner%null%1,                            Steven D%Stellman%null%1,                         Hannah T%Jordan%null%1,                         Robert M%Brackbill%null%1,                         James E%Cone%null%1,                         Indira%Debchoudhury%null%1,                         Mark R%Farfel%null%1,                         Carolyn M%Greene%null%1,                         James L%Hadler%null%1,                         Joseph%Kennedy%null%1,                         Jiehui%Li%null%1,                         Jonathan%Liff%null%1,                         Leslie%Stayner%null%1,                         Steven D%Stellman%null%1]"
$ws.Range("E22").Value = "[L. A.%Hendrickson%xref no email%1,       R. L.%Vogt%xref no email%1]"
$ws.Range("E23").Value = "[ Matthew N.%Peters%null%1,                            Morgan J.%Katz%null%1,                            John C.%Moscona%null%1,                            Mohi E.%Alkadri%null%1,                            Rashad H.%Khazi Syed%null%1,                            Thomas A.%Turnage%null%1,                            Vikram S.%Nijjar%null%1,                            Mohannad B.%Bisharat%null%1,                            Patrice%Delafontaine%null%0,                            Anand M.%Irimpen%null%1,                         Matthew N.%Peters%null%1,                         Morgan J.%Katz%null%1,                         John C.%Moscona%null%1,                         Mohi E.%Alkadri%null%1,                         Rashad H.%Khazi Syed%null%1,                         Thomas A.%Turnage%null%1,                         Vikram S.%Nijjar%null%1,                         Mohannad B.%Bisharat%null%1,                         Patrice%Delafontaine%null%0,                         Anand M.%Irimpen%null%1]"
$ws.Range("E24").Value = "[Ibrahim Baimba%Koroma%NULL%1,                           Dena%Javadi%NULL%2,                           Dena%Javadi%NULL%0,                           Katrina%Hann%NULL%5,                           Anthony D%Harries%NULL%2,                           Anthony D%Harries%NULL%0,                           Francis%Smart%NULL%2,                           Francis%Smart%NULL%0,                           Thomas%Samba%NULL%1,                           Jeffery%Edwards%NULL%4,                           Jeffery%Edwards%NULL%0,                           Priyakanta%Nayak%NULL%2,                           Priyakanta%Nayak%NULL%0,                           Katrina%Hann%NULL%0,                           Katrina%Hann%NULL%0,                           Jeffery%Edwards%NULL%0,                           Jeffery%Edwards%NULL%0,                           Katrina%Hann%NULL%0,                           Katrina%Hann%NULL%0]"
$ws.Range("E25").Value = "[Kyungeh%An%xref no email%1,       Jeanne%Salyer%xref no email%1,       Hsueh-Fen Sabrina%Kao%xref no email%1]"
$ws.Range("E26").Value = "[Shao%Lin%xref no email%1,       Marta I.%Gomez%xref no email%1,       Lenore%Gensburg%xref no email%1,       Wei%Liu%xref no email%1,       Syni-An%Hwang%xref no email%1]"
$ws.Range("E27").Value = "[Shengchao%Yu%xref no email%1,       Howard E.%Alper%xref no email%1,       Angela-Maithy%Nguyen%xref no email%1,       Robert M.%Brackbill%xref no email%1]"
$ws.Range("E28").Value = "[William%Gerin%xref no email%1,       William%Chaplin%xref no email%1,       Joseph E%Schwartz%xref no email%1,       John%Holland%xref no email%1,       Robert%Alter%xref no email%1,       Ronald%Wheeler%xref no email%1,       Diep%Duong%xref no email%1,       Thomas G%Pickering%xref no email%1]"
$ws.Range("E29").Value = "[ M. Natalie%Husarewycz%null%1,                            Renée%El-Gabalawy%null%1,                            Sarvesh%Logsetty%null%1,                            Jitender%Sareen%null%1,                         M. Natalie%Husarewycz%null%1,                         Renée%El-Gabalawy%null%1,                         Sarvesh%Logsetty%null%1,                         Jitender%Sareen%null%1]"
$ws.Range("E30").Value = "[J.%Ng%xref no email%1,       S. L.%Atkin%xref no email%1,       A. S.%Rigby%xref no email%1,       C.%Walton%xref no email%1,       E. S.%Kilpatrick%xref no email%1]"
$ws.Range("E31").Value = "[ Grégoire%Rey%null%1,                  Eric%Jougla%null%1,                  Anne%Fouillet%null%1,                  Gérard%Pavillon%null%1,                  Pierre%Bessemoulin%null%1,                  Philippe%Frayssinet%null%1,                  Jacqueline%Clavel%null%1,                  Denis%Hémon%null%1]"
$ws.Range("E32").Value = "[Ahmad%Jabbar%coreGivesNoEmail%1,                          Alaa%Boulad%coreGivesNoEmail%1,                          Anand%Irimpen%coreGivesNoEmail%1,                          Bradley%Deere%coreGivesNoEmail%1,                          Christopher%Westley%coreGivesNoEmail%1,                          Hassan%Baydoun%coreGivesNoEmail%1,                          Holly%Gonzales%coreGivesNoEmail%1,                          Indrajeet%Mahata%coreGivesNoEmail%1,                          John%C. Moscona%coreGivesNoEmail%1,                          Kapil%Yadav%coreGivesNoEmail%1,                          Matthew%N. Peters%coreGivesNoEmail%1,                          Patrick%Ters%coreGivesNoEmail%1,                          Paul%Katigbak%coreGivesNoEmail%1,                          Rohit%Maini%coreGivesNoEmail%1,                          Ryan%Nelson%coreGivesNoEmail%1,                          Sudesh%Srivastav%coreGivesNoEmail%1,                          Taraka%V. Gadiraju%coreGivesNoEmail%1]"
$ws.Range("E33").Value = "[Tina K.%Thethi%xref no email%1,       C. Lillian%Yau%xref no email%1,       Lizheng%Shi%xref no email%1,       Sharice%Leger%xref no email%1,       Prathima%Nagireddy%xref no email%1,       Jay%Waddadar%xref no email%1,       Prasanth%Surampudi%xref no email%1,       Jennifer%John-Kalarickal%xref no email%1,       Lydia%Yenoby%xref no email%1,       Vivian%Fonseca%xref no email%1]"
